# Apply edits described in the commit:
# "Fixed figure 20 to match SAS output- doesn't match figure from
#  Equity_Report_Figures_20231102.xlsx- need to note."
#
# 1) Education sheet: A1 header renamed from "edu" to "education_level"
# 2) Figure_20 sheet: numeric data values updated

$wb = $excel.ActiveWorkbook

# --- 1) Education sheet header rename ---
$wsEdu = $wb.Worksheets.Item("Education")
$wsEdu.Range("A1").Value = "education_level"

# --- 2) Figure_20 sheet data updates ---
$wsFig20 = $wb.Worksheets.Item("Figure_20")

# Row 2 - Imperial
$wsFig20.Range("D2").Value = 0.1
$wsFig20.Range("F2").Value = 2.26
$wsFig20.Range("H2").Value = 0.1

# Row 3 - Los Angeles
$wsFig20.Range("B3").Value = 0.22
$wsFig20.Range("C3").Value = 0.53
$wsFig20.Range("D3").Value = 0.29
$wsFig20.Range("E3").Value = 0.37
$wsFig20.Range("F3").Value = 0.13
$wsFig20.Range("G3").Value = 0.21
$wsFig20.Range("H3").Value = 0.28

# Row 4 - Orange
$wsFig20.Range("B4").Value = 0.11
$wsFig20.Range("C4").Value = 0
$wsFig20.Range("D4").Value = 0.14
$wsFig20.Range("E4").Value = 0.41
$wsFig20.Range("G4").Value = 0.08
$wsFig20.Range("H4").Value = 0.11

# Row 5 - Riverside
$wsFig20.Range("B5").Value = 0.12
$wsFig20.Range("D5").Value = 0.2
$wsFig20.Range("E5").Value = 0.3
$wsFig20.Range("F5").Value = 1.79
$wsFig20.Range("G5").Value = 0.18
$wsFig20.Range("H5").Value = 0.19

# Row 6 - San Bernardino
$wsFig20.Range("B6").Value = 0.31
$wsFig20.Range("C6").Value = 0.26
$wsFig20.Range("D6").Value = 0.18
$wsFig20.Range("E6").Value = 0.14
$wsFig20.Range("F6").Value = 0.53
$wsFig20.Range("G6").Value = 0.27
$wsFig20.Range("H6").Value = 0.23

# Row 7 - Ventura
$wsFig20.Range("D7").Value = 0.18
$wsFig20.Range("E7").Value = 0.33
$wsFig20.Range("G7").Value = 0.26
$wsFig20.Range("H7").Value = 0.22

# Row 8 - SCAG
$wsFig20.Range("B8").Value = 0.19
$wsFig20.Range("C8").Value = 0.42
$wsFig20.Range("D8").Value = 0.24
$wsFig20.Range("E8").Value = 0.34
$wsFig20.Range("F8").Value = 0.54
$wsFig20.Range("G8").Value = 0.19
$wsFig20.Range("H8").Value = 0.23
